$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-77 down to 67-78
$ws.Rows.Item(66).Insert()

$ws.Cells.Item(66, 1).Value = 11
$ws.Cells.Item(66, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(66, 3).Value = "Bíobío"
$ws.Cells.Item(66, 4).Value = 45258
$ws.Cells.Item(66, 5).Value = 8
$ws.Cells.Item(66, 6).Value = 100112026
$ws.Cells.Item(66, 7).Value = "Haba"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 100
$ws.Cells.Item(66, 11).Value = 8000
$ws.Cells.Item(66, 12).Value = 8000
$ws.Cells.Item(66, 13).Value = 8000
$ws.Cells.Item(66, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(66, 15).Value = "Región del Maule"
$ws.Cells.Item(66, 16).Value = 320
$ws.Cells.Item(66, 17).Value = 25
$ws.Cells.Item(66, 18).Value = "Hortaliza"
